# paises.xlsx -- "Update countries & provincias Spain"
# Refresh the COVID-19 country statistics snapshot and its timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: snapshot timestamp ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 4 de Octubre de 2020 a las 20:12"

# --- Per-country figure refresh (country/ranking unchanged) -----------------
# Row 4: Estados Unidos
$ws.Range("B4").Value = 7618026
$ws.Range("C4").Value = 16249
$ws.Range("D4").Value = 4833697
$ws.Range("E4").Value = 2569905
$ws.Range("G4").Value = 147
$ws.Range("H4").Value = 214424

# Row 27: Israel
$ws.Range("B27").Value = 265932
$ws.Range("C27").Value = 1489
$ws.Range("D27").Value = 194976
$ws.Range("E27").Value = 69249
$ws.Range("G27").Value = 25
$ws.Range("H27").Value = 1707

# Row 34: Marruecos
$ws.Range("B34").Value = 133272
$ws.Range("C34").Value = 2044
$ws.Range("D34").Value = 111036
$ws.Range("E34").Value = 19906
$ws.Range("G34").Value = 37
$ws.Range("H34").Value = 2330

# Row 64: Argelia
$ws.Range("B64").Value = 52136
$ws.Range("C64").Value = 141
$ws.Range("D64").Value = 36578
$ws.Range("E64").Value = 13798
$ws.Range("G64").Value = 4
$ws.Range("H64").Value = 1760

# Row 74: Irlanda
$ws.Range("B74").Value = 38032
$ws.Range("C74").Value = 364
$ws.Range("E74").Value = 12858

# Row 86: Costa de Marfil
$ws.Range("B86").Value = 19882
$ws.Range("C86").Value = 33
$ws.Range("D86").Value = 19449
$ws.Range("E86").Value = 313

# Row 105: Maldivas
$ws.Range("B105").Value = 10530
$ws.Range("C105").Value = 65
$ws.Range("D105").Value = 9364
$ws.Range("E105").Value = 1132

# Row 110: Haiti
$ws.Range("B110").Value = 8819
$ws.Range("C110").Value = 8
$ws.Range("D110").Value = 6992
$ws.Range("E110").Value = 1598

# Row 120: Malaui
$ws.Range("B120").Value = 5786
$ws.Range("C120").Value = 3
$ws.Range("D120").Value = 4541
$ws.Range("E120").Value = 1066

# Row 121: Suazilandia
$ws.Range("B121").Value = 5569
$ws.Range("C121").Value = 39
$ws.Range("D121").Value = 5127
$ws.Range("E121").Value = 331

# Row 123: Republica de Yibuti
$ws.Range("B123").Value = 5419
$ws.Range("C123").Value = 1
$ws.Range("D123").Value = 5349
$ws.Range("E123").Value = 9

# Row 141: Gambia
$ws.Range("B141").Value = 3594
$ws.Range("C141").Value = 4
$ws.Range("E141").Value = 1253

# Row 152: Sierra Leona
$ws.Range("B152").Value = 2269
$ws.Range("C152").Value = 10
$ws.Range("D152").Value = 1706
$ws.Range("E152").Value = 491

# Row 158: Yemen
$ws.Range("E158").Value = 130
$ws.Range("G158").Value = 2
$ws.Range("H158").Value = 591

# Row 165: Liberia
$ws.Range("B165").Value = 1348
$ws.Range("C165").Value = 1
$ws.Range("D165").Value = 1236
$ws.Range("E165").Value = 30

# Row 166: Republica del Chad
$ws.Range("B166").Value = 1217
$ws.Range("C166").Value = 3
$ws.Range("E166").Value = 56
$ws.Range("G166").Value = 1
$ws.Range("H166").Value = 86

# --- Rows 139-140: Somalia overtakes Estonia in total cases -> rows swap ---
$ws.Range("A139").Value = "Somalia"
$ws.Range("B139").Value = 3745
$ws.Range("C139").Value = 152
$ws.Range("D139").Value = 3010
$ws.Range("E139").Value = 636
$ws.Range("H139").Value = 99

$ws.Range("A140").Value = "Estonia"
$ws.Range("B140").Value = 3607
$ws.Range("C140").Value = 30
$ws.Range("D140").Value = 2749
$ws.Range("E140").Value = 791
$ws.Range("H140").Value = 67

# --- Rows 215-216: Montserrat overtakes Islas Malvinas -> rows swap --------
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
